$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-range labels in column B (bottom-up; B6 "0 to 5 minutes"
# stays as-is). The GAS minute-scale boundaries shift by one minute so the
# bands no longer overlap (5/15/25/35 -> 6/16/26/36).
$ws.Range("B5").Value = "6 to 15 minutes"
$ws.Range("B4").Value = "16 to 25 minutes"
$ws.Range("B3").Value = "26 to 35 minutes"
$ws.Range("B2").Value = "36 to 45 minutes"

# Update the selection to match the committed workbook state.
$ws.Range("B3").Select()
